# Updated lists of techs and RES
# Insert 6 new fuel rows (transportation techs) above the last row
# ("Wind energy" / "WND"), pushing it from row 39 down to row 45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert six blank rows starting at row 39; existing row 39 ("Wind energy")
# and everything below it shifts down by six rows (to row 45).
$ws.Range("A39:B44").Insert(-4121)

# Populate the new rows in the same order the values were originally
# entered (names for the two rail rows were typed together, then both
# codes filled in afterwards), so new shared-string entries line up.
$ws.Cells.Item(39, 1).Value = "All transportation by bus"
$ws.Cells.Item(39, 2).Value = "TRA_BUS"

$ws.Cells.Item(40, 1).Value = "All transportation by car"
$ws.Cells.Item(40, 2).Value = "TRA_CAR"

$ws.Cells.Item(41, 1).Value = "All transportation by motorcycle"
$ws.Cells.Item(41, 2).Value = "TRA_MCY"

$ws.Cells.Item(42, 1).Value = "All railway passenger transportation"
$ws.Cells.Item(43, 1).Value = "All railway freight transportation"
$ws.Cells.Item(42, 2).Value = "TRA_TRAIN_PSNG"
$ws.Cells.Item(43, 2).Value = "TRA_TRAIN_FREIGHT"

$ws.Cells.Item(44, 1).Value = "All truck transportation"
$ws.Cells.Item(44, 2).Value = "TRA_TRUCK"

# Match the author's final cursor position / selection in the diff.
$ws.Range("B45").Select()
